$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "미분방정식과 라플라스 변환"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/07/15/ODE_and_Laplace_transform.html"

$ws.Range("D36").Value = "Value-based Learning"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/325"

$ws.Range("D37").Value = "[Paper Review] End-to-End Object Detection with Transformers"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1784&mod=document&pageid=1"

$ws.Range("D51").Value = "[리눅스] root 패스워드 변경하는 법"
$ws.Range("E51").Value = "https://bskyvision.com/1213"
